# Reorder the comma-separated "Recorded By" values in column G by rotating
# the list left by one position (the first entry moves to the end).
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $text = [string]$val
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ", "
    if ($parts.Count -gt 1) {
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
        $cell.Value2 = $rotated
    }
}
